# Update column G ("K") values for rows 2-22 on the active sheet.
# These values represent a regenerated "K" (strikeouts) statistic
# replacing the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 5
    4  = 2
    5  = 4
    6  = 6
    7  = 2
    8  = 7
    9  = 5
    10 = 6
    11 = 6
    12 = 9
    13 = 8
    14 = 8
    15 = 12
    16 = 3
    17 = 6
    18 = 2
    19 = 8
    20 = 2
    21 = 4
    22 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
